# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (E16:E46) block is re-sorted from descending
# (2105 .. 1811) to ascending (1811 .. 2105) order, the monthly
# "Valor Mora" (F) figures for the first/last rows swap, and the
# "Salario Basico" (G) column is updated to the new salary value for
# every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 16..46 (was descending 2105 -> 1811)
$periods = @(
    "1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105"
)

$firstRow = 16
$lastRow = 46
$newSalario = 781242

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i

    # Column E: "Periodo Mora" label
    $ws.Cells.Item($row, 5).Value = $periods[$i]

    # Column G: "Salario Basico" -> new value for every row
    $ws.Cells.Item($row, 7).Value = $newSalario
}

# Column F: "Valor Mora" for the first and last row swap values
$ws.Cells.Item($firstRow, 6).Value = 29960
$ws.Cells.Item($lastRow, 6).Value = 24999
